$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothing = @{
    2 = 'Jumpsuit,Blouse'
    3 = 'Halter,Blazer'
    4 = 'Cutoffs,Blazer'
    5 = 'Cutoffs,Top'
    6 = 'Kaftan,Jumpsuit'
    7 = 'Blazer,Halter'
    8 = 'Halter,Blazer'
    9 = 'Jumpsuit,Dress'
    10 = 'Jumpsuit,Halter'
    11 = 'Jumpsuit,Blouse'
    12 = 'Jumpsuit,Blouse'
    13 = 'Tee,Jumpsuit'
    14 = 'Jumpsuit,Tee'
    15 = 'Caftan,Trunks'
    16 = 'Caftan,Halter'
    17 = 'Jumpsuit,Halter'
    18 = 'Caftan,Blouse'
    19 = 'Trunks,Caftan'
    20 = 'Trunks,Blazer'
    21 = 'Jumpsuit,Parka'
    22 = 'Trunks,Jumpsuit'
    23 = 'Halter,Top'
    24 = 'Blazer,Halter'
    25 = 'Halter,Jumpsuit'
    26 = 'Jumpsuit,Kaftan'
    27 = 'Tee,Blouse'
    28 = 'Trunks,Dress'
    29 = 'Halter,Jumpsuit'
    30 = 'Jumpsuit,Kaftan'
    31 = 'Jumpsuit,Blouse'
    32 = 'Tee,Blouse'
    33 = 'Jumpsuit,Kaftan'
    34 = 'Halter,Jumpsuit'
    35 = 'Caftan,Jumpsuit'
    36 = 'Halter,Blazer'
    37 = 'Trunks,Caftan'
    38 = 'Blouse,Caftan'
    39 = 'Jumpsuit,Kaftan'
    40 = 'Coat,Jumpsuit'
    41 = 'Coat,Halter'
    42 = 'Jumpsuit,Caftan'
    43 = 'Jumpsuit,Halter'
    44 = 'Jumpsuit,Blouse'
    45 = 'Jumpsuit,Kaftan'
    46 = 'Caftan,Trunks'
    47 = 'Jumpsuit,Halter'
    48 = 'Blazer,Trunks'
    49 = 'Caftan,Trunks'
    50 = 'Trunks,Caftan'
    51 = 'Halter,Cutoffs'
    52 = 'Blazer,Trunks'
    53 = 'Trunks,Caftan'
    54 = 'Caftan,Trunks'
    55 = 'Jumpsuit,Trunks'
    56 = 'Trunks,Jumpsuit'
    57 = 'Tee,Top'
    58 = 'Trunks,Jumpsuit'
    59 = 'Jumpsuit,Blouse'
    60 = 'Halter,Top'
    61 = 'Blazer,Blouse'
    62 = 'Halter,Blazer'
    63 = 'Caftan,Trunks'
    64 = 'Caftan,Jumpsuit'
    65 = 'Turtleneck,Blouse'
    66 = 'Turtleneck,Cutoffs'
    67 = 'Blazer,Top'
    68 = 'Blazer,Halter'
    69 = 'Parka,Blouse'
    70 = 'Turtleneck,Jumpsuit'
    71 = 'Jumpsuit,Kaftan'
    72 = 'Halter,Jumpsuit'
    73 = 'Halter,Blouse'
    74 = 'Trunks,Coverup'
    75 = 'Jumpsuit,Halter'
    76 = 'Halter,Jumpsuit'
    77 = 'Jumpsuit,Dress'
    78 = 'Jumpsuit,Trunks'
    79 = 'Caftan,Trunks'
    80 = 'Halter,Parka'
    81 = 'Jumpsuit,Halter'
    82 = 'Caftan,Trunks'
    83 = 'Jumpsuit,Halter'
    84 = 'Halter,Tee'
    85 = 'Jumpsuit,Tank'
    86 = 'Parka,Jumpsuit'
    87 = 'Jumpsuit,Blazer'
    88 = 'Trunks,Caftan'
    89 = 'Jumpsuit,Top'
    90 = 'Jumpsuit,Tee'
    91 = 'Blouse,Halter'
    92 = 'Coat,Cutoffs'
    93 = 'Jumpsuit,Kaftan'
    94 = 'Jumpsuit,Halter'
    95 = 'Trunks,Coverup'
    96 = 'Jumpsuit,Kaftan'
    97 = 'Jumpsuit,Dress'
    98 = 'Caftan,Jumpsuit'
    99 = 'Trunks,Coverup'
    100 = 'Turtleneck,Capris'
    101 = 'Jumpsuit,Blouse'
    102 = 'Jumpsuit,Kaftan'
    103 = 'Cutoffs,Trunks'
    104 = 'Trunks,Cutoffs'
    105 = 'Kaftan,Blouse'
    106 = 'Caftan,Jumpsuit'
    107 = 'Jumpsuit,Trunks'
    108 = 'Top,Halter'
    109 = 'Top,Caftan'
    110 = 'Top,Parka'
    111 = 'Parka,Sweatpants'
    112 = 'Kaftan,Dress'
    113 = 'Cutoffs,Coat'
    114 = 'Caftan,Trunks'
    115 = 'Trunks,Cutoffs'
    116 = 'Halter,Blazer'
    117 = 'Halter,Jumpsuit'
    118 = 'Turtleneck,Halter'
    119 = 'Turtleneck,Halter'
    120 = 'Gauchos,Blouse'
    121 = 'Jumpsuit,Trunks'
    122 = 'Halter,Top'
    123 = 'Top,Leggings'
    124 = 'Halter,Turtleneck'
    125 = 'Jumpsuit,Caftan'
    126 = 'Caftan,Halter'
    127 = 'Trunks,Blazer'
    128 = 'Caftan,Trunks'
    129 = 'Jumpsuit,Halter'
    130 = 'Halter,Blazer'
    131 = 'Caftan,Jumpsuit'
    132 = 'Jumpsuit,Halter'
    133 = 'Jumpsuit,Blouse'
    134 = 'Halter,Blazer'
    135 = 'Blazer,Hoodie'
    136 = 'Trunks,Jumpsuit'
    137 = 'Jumpsuit,Caftan'
    138 = 'Jumpsuit,Blouse'
    139 = 'Jumpsuit,Blouse'
    140 = 'Blouse,Caftan'
    141 = 'Jumpsuit,Trunks'
    142 = 'Caftan,Parka'
    143 = 'Caftan,Parka'
    144 = 'Caftan,Trunks'
    145 = 'Caftan,Jumpsuit'
    146 = 'Jumpsuit,Dress'
    147 = 'Jumpsuit,Jodhpurs'
    148 = 'Parka,Caftan'
    149 = 'Kaftan,Jumpsuit'
    150 = 'Caftan,Trunks'
    151 = 'Halter,Top'
    152 = 'Caftan,Turtleneck'
    153 = 'Halter,Blouse'
    154 = 'Jumpsuit,Tee'
    155 = 'Dress,Trunks'
    156 = 'Jumpsuit,Tee'
    157 = 'Blouse,Caftan'
    158 = 'Top,Coat'
    159 = 'Sweatpants,Gauchos'
    160 = 'Halter,Jumpsuit'
    161 = 'Coat,Jumpsuit'
    162 = 'Jumpsuit,Tee'
    163 = 'Jumpsuit,Halter'
    164 = 'Caftan,Trunks'
    165 = 'Parka,Gauchos'
    166 = 'Gauchos,Blouse'
    167 = 'Trunks,Coverup'
    168 = 'Halter,Blazer'
    169 = 'Jumpsuit,Halter'
    170 = 'Cutoffs,Coat'
    171 = 'Blouse,Dress'
    172 = 'Jumpsuit,Kaftan'
    173 = 'Jumpsuit,Kaftan'
    174 = 'Halter,Blazer'
    175 = 'Blazer,Halter'
}

foreach ($row in $clothing.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $clothing[$row]
}
